$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing it to stay a text cell (matches the
# source workbook, which stores every data value as inlineStr) and without
# leaving any residual number-format/style change behind.
function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "246.23"

Set-TextValue $ws.Range("D3") "21.63"

Set-TextValue $ws.Range("D4") "5.275"

Set-TextValue $ws.Range("D5") "0.05590"

Set-TextValue $ws.Range("D6") "3.388"

Set-TextValue $ws.Range("D7") "6.382"

Set-TextValue $ws.Range("D8") "0.8138"

Set-TextValue $ws.Range("D9") "0.9721"

$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D10") "0.01137"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1409"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.07376"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D13") "0.03146"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03035"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09289"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D16") "3.561"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D17") "0.001610"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D18") "0.04711"
$ws.Range("E18").Value = "17CoinExTokenCET"

Set-TextValue $ws.Range("D19") "0.006366"

Set-TextValue $ws.Range("D20") "0.005060"

Set-TextValue $ws.Range("D21") "0.001030"

Set-TextValue $ws.Range("D22") "0.0001495"

Set-TextValue $ws.Range("D23") "3.773"

Set-TextValue $ws.Range("D24") "2.123"

Set-TextValue $ws.Range("D25") "0.3251"

Set-TextValue $ws.Range("D26") "0.1299"

Set-TextValue $ws.Range("D28") "0.0003090"

Set-TextValue $ws.Range("D40") "0.03914"

Set-TextValue $ws.Range("D41") "0.007024"

Set-TextValue $ws.Range("D43") "0.003007"

Set-TextValue $ws.Range("D44") "0.007752"

Set-TextValue $ws.Range("D45") "0.00005788"

Set-TextValue $ws.Range("D46") "0.00000000748"

Set-TextValue $ws.Range("D47") "0.0005484"

Set-TextValue $ws.Range("D48") "0.6778"

Set-TextValue $ws.Range("D49") "0.1387"
$ws.Range("E49").Value = "48BOLOBOLO"

Set-TextValue $ws.Range("D50") "0.00002093"

Set-TextValue $ws.Range("D51") "0.01007"
